$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 11419
$ws.Range("F3").Value = 1940
$ws.Range("F7").Value = 766
$ws.Range("F8").Value = 1011
$ws.Range("F9").Value = 599
$ws.Range("F11").Value = 1348
$ws.Range("F12").Value = 679
$ws.Range("F15").Value = 991
$ws.Range("F16").Value = 539
$ws.Range("F18").Value = 1128
$ws.Range("F19").Value = 213
$ws.Range("F20").Value = 936
$ws.Range("F22").Value = 138
$ws.Range("F23").Value = 306
$ws.Range("F25").Value = 266
$ws.Range("F26").Value = 471
$ws.Range("F27").Value = 505
$ws.Range("F28").Value = 678
$ws.Range("F29").Value = 177
$ws.Range("F30").Value = 112
$ws.Range("F31").Value = 327
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 897
$ws.Range("F7").Value = 7
$ws.Range("F8").Value = 99
$ws.Range("F9").Value = 43
$ws.Range("F10").Value = 390
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 50
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 11419
$ws.Range("F3").Value = 1940
$ws.Range("F8").Value = 766
$ws.Range("F9").Value = 1011
$ws.Range("F11").Value = 599
$ws.Range("F13").Value = 50
$ws.Range("F14").Value = 1348
$ws.Range("F16").Value = 679
$ws.Range("F18").Value = 897
$ws.Range("F20").Value = 991
$ws.Range("F21").Value = 539
$ws.Range("F23").Value = 1128
$ws.Range("F24").Value = 213
$ws.Range("F25").Value = 936
$ws.Range("F27").Value = 138
$ws.Range("F28").Value = 306
$ws.Range("F31").Value = 266
$ws.Range("F32").Value = 7
$ws.Range("F33").Value = 99
$ws.Range("F34").Value = 99
$ws.Range("F35").Value = 471
$ws.Range("F36").Value = 505
$ws.Range("F37").Value = 678
$ws.Range("F38").Value = 177
$ws.Range("F39").Value = 43
$ws.Range("F40").Value = 112
$ws.Range("F41").Value = 390
$ws.Range("F42").Value = 327
